$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header cells for the "Synoptic" weekly-comment columns (R1, S1)
# ---------------------------------------------------------------------------
$ws.Range("R1").Value = "Synoptic Wk1"
$ws.Range("S1").Value = "Synoptic Wk2"

# ---------------------------------------------------------------------------
# 2. New values in the existing "Application Development" block (rows 17-19)
#    and "Deliverables" block (row 21), plus their highlight colours.
#    Formatting is carried over from existing coloured cells on the same
#    rows via copy / paste-special (formats only) so the new cells reuse the
#    same style records as the rest of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("N16").Copy() | Out-Null
$ws.Range("R17").PasteSpecial(-4122) | Out-Null

$ws.Range("R18").Value = "Log File Addition"
$ws.Range("N18").Copy() | Out-Null
$ws.Range("R18").PasteSpecial(-4122) | Out-Null

$ws.Range("R21").Value = "Editing of Intro & Concl"
$ws.Range("S21").Value = "Tweaking of sections"
$ws.Range("P21").Copy() | Out-Null
$ws.Range("R21").PasteSpecial(-4122) | Out-Null
$ws.Range("P21").Copy() | Out-Null
$ws.Range("S21").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Threaded comments recording what changed each week.
# ---------------------------------------------------------------------------
$ws.Range("R18").AddCommentThreaded("Tried adding log files. Location in paper: Limitation") | Out-Null
$ws.Range("R21").AddCommentThreaded("Changed Intro questions, planning hypothesis, added limitations and improvements.") | Out-Null
$ws.Range("S21").AddCommentThreaded("Mild rewording of conclusion, added extenstions, more text in findings") | Out-Null

# ---------------------------------------------------------------------------
# 4. New "Extra Comments" section at the bottom of the sheet (rows 23-28).
# ---------------------------------------------------------------------------
$ws.Range("A23").Value = "Extra Comments"

$ws.Range("A24").Value = "Log Files (XML/CSV)"
$ws.Range("B24").Value = "Could not be implemented due to outdated codes found"

$ws.Range("A25").Value = "Intro"
$ws.Range("B25").Value = "Changed questions, added hypothesis."

$ws.Range("A26").Value = "Conclusion"
$ws.Range("B26").Value = "Added extention, limitation, improvement. Reworded slightly."

$ws.Range("A27").Value = "Citations"
$ws.Range("B27").Value = "Added intext citations instead of leaving only numbers."

$ws.Range("A28").Value = "Findings"
$ws.Range("B28").Value = "Added few more information that was collected."

# ---------------------------------------------------------------------------
# 5. Column widths: widen column B to fit the longer text, and size the two
#    new "Synoptic" columns.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 51.21875
$ws.Columns.Item(18).ColumnWidth = 20.109375
$ws.Columns.Item(19).ColumnWidth = 18.21875

# ---------------------------------------------------------------------------
# 6. Selection moved to B20 to match the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("B20").Select() | Out-Null
